# "updated chaining of login"
#
# The News headline text used by the NewsSearch test page ("news headlines
# sports news") is retired from the shared-string pool and the AddNewNews
# test page's expected headline is flipped from "...started" to
# "...completed" now that the login chaining covers the whole flow.
#
# We reproduce this by writing the NewsSearch page's headline cell back to
# its own (unchanged) text -- which forces the workbook to rebuild the
# shared-string table -- and by updating the AddNewNews page's headline to
# the new "completed" wording. Excel's automatic dedup/compaction of the
# shared-string table then naturally drops the now-unused old string and
# appends the new one at the end, matching the recorded diff.
$wb = $excel.ActiveWorkbook

$newsSearch = $wb.Worksheets.Item("NewsSearch")
$newsSearch.Range("A2").Value = "Obsqura automation Project started"

$addNewNews = $wb.Worksheets.Item("AddNewNews")
$addNewNews.Range("A2").Value = "Obsqura automation Project completed"

# The active/selected worksheet moves from AddNewSubCategory to AddNewNews.
$addNewNews.Activate() | Out-Null
$addNewNews.Range("A2").Select() | Out-Null
